$wb = $excel.ActiveWorkbook

# --- consumptionAssets: update logistics fleet yearly electricity demand formula ---
$wsConsumption = $wb.Worksheets.Item("consumptionAssets")
$wsConsumption.Range("F10").Formula = "=60*100000*1.3"

# --- storageAssets: add new 10MWh grid battery asset row ---
$wsStorage = $wb.Worksheets.Item("storageAssets")
$wsStorage.Range("A16").Value = 13
$wsStorage.Range("B16").Value = "Grid_battery_10MWh"
$wsStorage.Range("C16").Value = "STORAGE"
$wsStorage.Range("D16").Value = "STORAGE_ELECTRIC"
$wsStorage.Range("E16").Value = 2000
$wsStorage.Range("F16").Value = 0
$wsStorage.Range("G16").Value = 1
$wsStorage.Range("H16").Value = 0
$wsStorage.Range("I16").Value = 0
$wsStorage.Range("J16").Value = 0
$wsStorage.Range("K16").Value = 0
$wsStorage.Range("L16").Value = 10000
$wsStorage.Range("M16").Value = 0
$wsStorage.Range("M16").NumberFormat = "0.00E+00"

# --- Update the view state: selection moves from consumptionAssets!F10 to F11,
#     and the storageAssets sheet becomes the active tab, scrolled so column C
#     is leftmost, with L17 selected ---
$wsConsumption.Range("F11").Select()

$wsStorage.Activate()
$wsStorage.Range("C1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$wsStorage.Range("L17").Select()
